# Opinionated conversion of value type edge cases (#637)
# The "GeneratingUnit" view in the Properties sheet gains three new
# properties that the updated value-type resolution logic now surfaces:
# cable, operator and creationTime (each a plain "text" value type,
# not nullable, not a list) - mirroring the existing rows for that view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")

# Carry over the formatting (fill + border) used by the existing data
# rows so the three new rows look the same as rows 3-6.
$ws.Range("A6:P6").Copy()
$ws.Range("A7:P9").PasteSpecial(-4122)   # xlPasteFormats

# Row 7: GeneratingUnit.cable
$ws.Range("A7").Value = "GeneratingUnit"
$ws.Range("B7").Value = "cable"
$ws.Range("F7").Value = "text"
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $false
$ws.Range("K7").Value = "GeneratingUnit"
$ws.Range("L7").Value = "cable"
$ws.Range("M7").Value = "GeneratingUnit"
$ws.Range("N7").Value = "cable"

# Row 8: GeneratingUnit.operator
$ws.Range("A8").Value = "GeneratingUnit"
$ws.Range("B8").Value = "operator"
$ws.Range("F8").Value = "text"
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false
$ws.Range("K8").Value = "GeneratingUnit"
$ws.Range("L8").Value = "operator"
$ws.Range("M8").Value = "GeneratingUnit"
$ws.Range("N8").Value = "operator"

# Row 9: GeneratingUnit.creationTime
$ws.Range("A9").Value = "GeneratingUnit"
$ws.Range("B9").Value = "creationTime"
$ws.Range("F9").Value = "text"
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = $false
$ws.Range("K9").Value = "GeneratingUnit"
$ws.Range("L9").Value = "creationTime"
$ws.Range("M9").Value = "GeneratingUnit"
$ws.Range("N9").Value = "creationTime"

# The Properties sheet is now the one the author was looking at when the
# workbook was last saved: it becomes the active tab (was "Containers"),
# and the selection sits on the newly added data (L10).
$ws.Activate()
$ws.Range("L10").Select() | Out-Null
